$wb = $excel.ActiveWorkbook

# Rename sheets: "Master Material" -> "Master Item", "Master Vendor x Material" -> "Master Vendor x Item".
# Excel auto-updates formulas (VLOOKUP references) that point at the renamed
# sheets as part of the rename operation.
$wb.Worksheets.Item("Master Material").Name = "Master Item"
$wb.Worksheets.Item("Master Vendor x Material").Name = "Master Vendor x Item"

# Move the active selection on the "Master Vendor x Item" sheet (previously
# named "Master Vendor x Material") from C4 to F23.
$ws = $wb.Worksheets.Item("Master Vendor x Item")
$ws.Activate()
$ws.Range("F23").Select()
